$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.241.20'
$ws.Range('E2').Value = '  +0.22%  '
$ws.Range('D3').Value = '1.905.16'
$ws.Range('E3').Value = '  +0.20%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.26%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.53'
$ws.Range('E5').Value = '  -0.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.24%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5261'
$ws.Range('E7').Value = '  +1.14%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3818'
$ws.Range('E8').Value = '  +1.48%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07289'
$ws.Range('E9').Value = '  +0.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.02'
$ws.Range('E10').Value = '  +4.07%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9009'
$ws.Range('E11').Value = '  -0.33%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08195'
$ws.Range('E12').Value = '  -1.54%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '95.73'
$ws.Range('E13').Value = '  -1.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.349'
$ws.Range('E14').Value = '  +1.11%  '
$ws.Range('E15').Value = '  +0.20%  '
$ws.Range('B16').Value = 'Avalanche'
$ws.Range('C16').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.78'
$ws.Range('E16').Value = '  +1.69%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008630'
$ws.Range('E17').Value = '  -0.23%  '
$ws.Range('D18').Value = '1.351.46'
$ws.Range('E18').Value = '  -29.04%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.002'
$ws.Range('E19').Value = '  +0.22%  '
$ws.Range('D20').Value = '27.287.87'
$ws.Range('E20').Value = '  +0.25%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.075'
$ws.Range('E21').Value = '  -0.18%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.81'
$ws.Range('E22').Value = '  +1.61%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.521'
$ws.Range('E23').Value = '  +1.19%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '149.93'
$ws.Range('E24').Value = '  +2.38%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.298'
$ws.Range('E25').Value = '  -1.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '18.23'
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('E27').Value = '  -0.67%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '116.63'
$ws.Range('E28').Value = '  +1.35%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.833'
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.822'
$ws.Range('E30').Value = '  -1.46%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09256'
$ws.Range('E31').Value = '  +0.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.8318'
$ws.Range('E32').Value = '  +4.16%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05080'
$ws.Range('E33').Value = '  +0.14%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.225'
$ws.Range('E34').Value = '  -1.74%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.010'
$ws.Range('E35').Value = '  +2.32%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.354'
$ws.Range('E36').Value = '  -1.86%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.678'
$ws.Range('E37').Value = '  +3.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5800'
$ws.Range('E38').Value = '  +1.04%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02003'
$ws.Range('E39').Value = '  +0.15%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.079'
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '9.293'
$ws.Range('E41').Value = '  +3.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.533'
$ws.Range('E42').Value = '  -1.10%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '116.83'
$ws.Range('E43').Value = '  -0.21%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1521'
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4917'
$ws.Range('E45').Value = '  +1.11%  '
$ws.Range('E46').Value = '  +0.16%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.15'
$ws.Range('E47').Value = '  +0.18%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.635'
$ws.Range('E48').Value = '  +0.19%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '38.80'
$ws.Range('E49').Value = '  +2.85%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06198'
$ws.Range('E50').Value = '  +4.07%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '63.84'
$ws.Range('E51').Value = '  -0.29%  '

Write-Host "Updated cryptos list"
